$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 3
    10 = 2
    11 = 0
    12 = 1
    13 = 1
    14 = 3
    15 = 4
    16 = 0
    17 = 4
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 2
    24 = 1
    25 = 2
    26 = 2
    27 = 3
    28 = 3
    29 = 3
    30 = 2
    31 = 2
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    36 = 2
    37 = 0
    38 = 3
    39 = 0
    40 = 3
    41 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}

$wb.Save()
